# Applies the "enchant cost change" edit described in the commit message.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the CostDiamond (C column) values for rows 5-15
$ws.Range("C5").Value = 8000
$ws.Range("C6").Value = 8000
$ws.Range("C7").Value = 20000
$ws.Range("C8").Value = 24000
$ws.Range("C9").Value = 29000
$ws.Range("C10").Value = 35000
$ws.Range("C11").Value = 42000
$ws.Range("C12").Value = 50000
$ws.Range("C13").Value = 60000
$ws.Range("C14").Value = 72000
$ws.Range("C15").Value = 72000

# Update the view: zoom to 130% and move the active selection to D14
$ws.Activate()
$excel.ActiveWindow.Zoom = 130
$ws.Range("D14").Select()
